$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# Slide 11 ("Building Kernel") - TextBox 2
# -----------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$tb11 = $s11.Shapes.Item(2)

# 1. Split "If cross-compiling from Linux do not set environmental
#    variables on system level." into three runs by rewriting just the
#    "do " fragment as "does " (same rPr on every side, so the engine
#    only needs to re-split the run boundaries around the edited text).
$tr11 = $tb11.TextFrame.TextRange
$para3 = $tr11.Paragraphs(3)
$fragment = $tr11.Characters($para3.Start + 30, 3)
$fragment.Text = "does "

# 2. Widen the text box (cx grows from 8418051 to 8661836 EMU, cy stays at
#    5444267). The text box has spAutoFit, so the engine recomputes both
#    dimensions from the text layout as soon as either the text or the
#    size changes; re-assert both Width and Height (in points, 1 pt =
#    12700 EMU) to land back on the exact target EMU values.
$tb11.Width = 682.0343636220473
$tb11.Height = 428.6824784448819

# -----------------------------------------------------------------------
# Slide 30 ("Home work") - TextBox 2
# -----------------------------------------------------------------------
$s30 = $p.Slides.Item(30)
$tb30 = $s30.Shapes.Item(2)
$tr30 = $tb30.TextFrame.TextRange

# 3. Merge the "https://" + "wiki.openwrt.org/toh/tp-link/tl-mr3020" runs
#    into a single run (rewrite the whole paragraph span so the engine
#    collapses it back to one run, keeping the first run's formatting,
#    including the hyperlink).
$para2 = $tr30.Paragraphs(2)
$whole2 = $tr30.Characters($para2.Start, $para2.Length)
$whole2.Text = "https://wiki.openwrt.org/toh/tp-link/tl-mr3020"

# 4. Merge the "https://" + "github.com/duxing2007/ldd3-examples-3.x" runs
#    into a single run the same way.
$para5 = $tr30.Paragraphs(5)
$whole5 = $tr30.Characters($para5.Start, $para5.Length)
$whole5.Text = "https://github.com/duxing2007/ldd3-examples-3.x"

# 5. Merge the "'s " + "kernel drivers directory "" runs (inside paragraph
#    7) into a single run.
$para7 = $tr30.Paragraphs(7)
$mergeStart = $para7.Start + 42
$mergeLen = 29
$merged = $tr30.Characters($mergeStart, $mergeLen)
$merged.Text = [char]0x2019 + "s kernel drivers directory " + [char]0x201C

# This text box also has spAutoFit; none of the edits above change the
# visible character count (pure run-merges), but the engine still
# recomputes the box height from its own text-layout pass whenever the
# TextRange is touched. Re-assert the original height (in points) so the
# shape's on-disk size is left exactly as it was.
$tb30.Height = 302.3608324015748
